# "Readmes added to templates"
#
# Insert a new "README" worksheet in front of the existing "Occurrences"
# sheet and fill it with the GBIF IPT template README content/formatting.

$wb  = $excel.ActiveWorkbook
$occ = $wb.Worksheets.Item(1)

$readme = $wb.Worksheets.Add($occ)
$readme.Name = "README"

# ---------------------------------------------------------------------------
# column widths
# ---------------------------------------------------------------------------
$readme.Columns.Item(1).ColumnWidth = 27.17
$readme.Columns.Item(2).ColumnWidth = 83

# ---------------------------------------------------------------------------
# colors / constants used below
# ---------------------------------------------------------------------------
$yellow = 10092543   # RGB(255,255,153) FFFFFF99
$gray   = 8421504    # RGB(128,128,128) FF808080
$black  = 0           # RGB(0,0,0)       FF000000

# ---------------------------------------------------------------------------
# Block 1 (rows 1-5, cols A:B) -- title block on yellow background
# ---------------------------------------------------------------------------
$readme.Range("A1:B5").Interior.Color = $yellow
$readme.Range("A1:B5").Font.Name = "Calibri"
$readme.Range("A1:B5").Font.Size = 12

$readme.Range("B1").Value = "GBIF IPT Template:"
$readme.Range("B1").Font.Bold = $true

$readme.Rows.Item(2).RowHeight = 20
$readme.Range("B2").Value = "Occurrence Data"
$readme.Range("B2").Font.Bold = $true
$readme.Range("B2").Font.Size = 16

$readme.Rows.Item(4).RowHeight = 45
$readme.Range("B4").Value = "Use this template for filling in occurrence data. Upload the template to the IPT where it can be published in Darwin Core Archive (DwC-A) format. Note this template must be mapped to the Occurrence Core in the IPT."
$readme.Range("B4").Font.Italic = $true
$readme.Range("B4").WrapText = $true

$readme.Range("B3").Font.Bold = $true
$readme.Range("B3").Font.Size = 10
$readme.Range("B3").Font.Color = $gray
$readme.Range("B3").HorizontalAlignment = -4152

$readme.Range("B5").Font.Italic = $true
$readme.Range("B5").Font.Size = 11
$readme.Range("B5").WrapText = $true

# outline border around the whole title block
$readme.Range("A1:B5").Borders.Item(7).LineStyle = 1
$readme.Range("A1:B5").Borders.Item(7).Weight = 2
$readme.Range("A1:B5").Borders.Item(10).LineStyle = 1
$readme.Range("A1:B5").Borders.Item(10).Weight = 2
$readme.Range("A1:B1").Borders.Item(8).LineStyle = 1
$readme.Range("A1:B1").Borders.Item(8).Weight = 2
$readme.Range("B2:B5").Borders.Item(7).LineStyle = 1
$readme.Range("B2:B5").Borders.Item(7).Weight = 2
$readme.Range("B2").Borders.Item(9).LineStyle = 1
$readme.Range("B2").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# Block 2 (rows 6-8, cols A:B) -- "Sheet Name / Purpose" mini table
# ---------------------------------------------------------------------------
$readme.Range("A6:B6").Interior.Color = $yellow
$readme.Rows.Item(6).RowHeight = 16

$readme.Range("A6").Value = "Sheet Name"
$readme.Range("B6").Value = "Purpose"
$readme.Range("A6:B6").Font.Name = "Calibri"
$readme.Range("A6:B6").Font.Bold = $true
$readme.Range("A6:B6").Font.Size = 12
$readme.Range("A6:B6").HorizontalAlignment = -4131

$readme.Rows.Item(7).RowHeight = 30
$readme.Range("A7").Value = "Occurrences"
$readme.Range("A7").Font.Name = "Calibri"
$readme.Range("A7").Font.Bold = $true
$readme.Range("A7").Font.Size = 12
$readme.Range("A7").VerticalAlignment = -4160

$readme.Range("B7").Value = "This sheet is used to record a list of species at a particular place and normally on a specified date. Ideally locality information includes coordinates to support mapping. "
$readme.Range("B7:B8").Font.Name = "Calibri"
$readme.Range("B7:B8").Font.Size = 12
$readme.Range("B7:B8").WrapText = $true

$readme.Range("A8").Font.Name = "Calibri"
$readme.Range("A8").Font.Bold = $true
$readme.Range("A8").Font.Size = 12
$readme.Range("A8").VerticalAlignment = -4160

# borders around the mini table + thick rule under the header row
$readme.Range("A6:B8").Borders.Item(7).LineStyle = 1
$readme.Range("A6:B8").Borders.Item(7).Weight = 2
$readme.Range("A6:B8").Borders.Item(10).LineStyle = 1
$readme.Range("A6:B8").Borders.Item(10).Weight = 2
$readme.Range("A6:B6").Borders.Item(8).LineStyle = 1
$readme.Range("A6:B6").Borders.Item(8).Weight = 2
$readme.Range("A6:B6").Borders.Item(9).LineStyle = 1
$readme.Range("A6:B6").Borders.Item(9).Weight = -4138

$readme.Range("A9:A13").Borders.Item(7).LineStyle = 1
$readme.Range("A9:A13").Borders.Item(7).Weight = 2
$readme.Range("A9:A13").Borders.Item(7).ColorIndex = 48

# ---------------------------------------------------------------------------
# Notes block (rows 9-13)
# ---------------------------------------------------------------------------
$readme.Range("B9").Value = "Notes:"
$readme.Range("B9").Font.Name = "Calibri"
$readme.Range("B9").Font.Bold = $true
$readme.Range("B9").Font.Size = 12

$readme.Rows.Item(10).RowHeight = 30
$readme.Range("B10").Value = "#1. The header row shows required and recommended terms. Hover over the cell to find out if it's required or recommended, and to obtain a definition of the term examples. "
$readme.Range("B10").Font.Name = "Calibri"
$readme.Range("B10").Font.Bold = $true
$readme.Range("B10").Font.Size = 12
$readme.Range("B10").Font.Color = $black
$readme.Range("B10").WrapText = $true
$readme.Range("B10").Characters(5, 300).Font.Bold = $false
$readme.Range("B10").Characters(5, 300).Font.Color = $black

$readme.Rows.Item(11).RowHeight = 30
$readme.Range("B11").Value = "#2. Additional columns can be added, but you should use DwC term names: http://rs.tdwg.org/dwc/terms/"
$readme.Range("B11").Font.Name = "Calibri"
$readme.Range("B11").Font.Bold = $true
$readme.Range("B11").Font.Size = 12
$readme.Range("B11").WrapText = $true
$readme.Range("B11").Characters(5, 300).Font.Bold = $false

$readme.Range("B12").Value = "#3. Columns can be reordered, but the header name (equal to a DwC term name) cannot be changed."
$readme.Range("B12").Font.Name = "Calibri"
$readme.Range("B12").Font.Bold = $true
$readme.Range("B12").Font.Size = 12
$readme.Range("B12").WrapText = $true
$readme.Range("B12").Characters(5, 300).Font.Bold = $false

$readme.Range("A10:A13").Font.Name = "Calibri"
$readme.Range("A10:A13").Font.Size = 12

$readme.Range("B13").Font.Name = "Calibri"
$readme.Range("B13").Font.Size = 12

# ---------------------------------------------------------------------------
# sheet view / selection, then leave README as the active/visible tab
# ---------------------------------------------------------------------------
$readme.Range("B13").Select()
$readme.Activate()
